$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Locations sheet: update "Improved Time" (column D) values and give
# them a dedicated 3-decimal number format (new style index 3).
# ---------------------------------------------------------------------
$wsLoc = $wb.Worksheets.Item("Locations")

$wsLoc.Range("D3").Value = 0.835
$wsLoc.Range("D4").Value = 2.33
$wsLoc.Range("D5").Value = 10.251
$wsLoc.Range("D6").Value = 20.479
$wsLoc.Range("D7").Value = 102.065
$wsLoc.Range("D8").Value = 204.263
$wsLoc.Range("D3:D8").NumberFormat = "0.000"

# Column widths for C and D (bestFit, custom width) on the Locations sheet.
$wsLoc.Columns.Item(3).ColumnWidth = 14
$wsLoc.Columns.Item(4).ColumnWidth = 12.85546875

# Move the active selection from D15 to D8.
$wsLoc.Range("D8").Select()

# ---------------------------------------------------------------------
# Rewards sheet: update "Improved Time" (column D) values and apply the
# same dedicated 3-decimal number format.
# ---------------------------------------------------------------------
$wsRew = $wb.Worksheets.Item("Rewards")

$wsRew.Range("D3").Value = 1.638
$wsRew.Range("D4").Value = 11.039
$wsRew.Range("D5").Value = 103.766
$wsRew.Range("D6").Value = 1022.355
$wsRew.Range("D3:D6").NumberFormat = "0.000"
